$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number (45182 -> 2023-09-13).
# The commit updates every data row (2 through 533) to serial 45184 (2023-09-15).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 533 }

$range = $ws.Range("C2:C$lastRow")
$range.Value = 45184
